$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Tgfb3"
$ws.Cells.Item(2,3).Value = "Tgfbr3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.6423476666666667
$ws.Cells.Item(2,8).Value = 1.927043
$ws.Cells.Item(2,9).Value = 0.01173234890143342
$ws.Cells.Item(2,10).Value = 0.01173234890143342
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 14.85604233333333
$ws.Cells.Item(2,14).Value = 44.568127
$ws.Cells.Item(2,15).Value = 0.09286934904108346
$ws.Cells.Item(2,16).Value = 0.09286934904108346
$ws.Cells.Item(2,17).Value = 9.542744128717889
$ws.Cells.Item(2,18).Value = 85.884697158461
$ws.Cells.Item(2,19).Value = 0.001089575605198992
$ws.Cells.Item(2,20).Value = 0.001089575605198992

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Tgfb3"
$ws.Cells.Item(3,3).Value = "Tgfbr3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.6423476666666667
$ws.Cells.Item(3,8).Value = 1.927043
$ws.Cells.Item(3,9).Value = 0.01173234890143342
$ws.Cells.Item(3,10).Value = 0.01173234890143342
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 114.2734143333333
$ws.Cells.Item(3,14).Value = 342.820243
$ws.Cells.Item(3,15).Value = 0.7143556381787382
$ws.Cells.Item(3,16).Value = 0.7143556381787382
$ws.Cells.Item(3,17).Value = 73.40326105904988
$ws.Cells.Item(3,18).Value = 660.6293495314491
$ws.Cells.Item(3,19).Value = 0.008381069586819088
$ws.Cells.Item(3,20).Value = 0.00838106958681909

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Tgfb3"
$ws.Cells.Item(4,3).Value = "Tgfbr3"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.6423476666666667
$ws.Cells.Item(4,8).Value = 1.927043
$ws.Cells.Item(4,9).Value = 0.01173234890143342
$ws.Cells.Item(4,10).Value = 0.01173234890143342
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 30.83766366666667
$ws.Cells.Item(4,14).Value = 92.512991
$ws.Cells.Item(4,15).Value = 0.1927750127801784
$ws.Cells.Item(4,16).Value = 0.1927750127801784
$ws.Cells.Item(4,17).Value = 19.80850130173478
$ws.Cells.Item(4,18).Value = 178.276511715613
$ws.Cells.Item(4,19).Value = 0.002261703709415339
$ws.Cells.Item(4,20).Value = 0.00226170370941534

# Row 5
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Tgfb3"
$ws.Cells.Item(5,3).Value = "Tgfbr3"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 22.27635266666667
$ws.Cells.Item(5,8).Value = 66.829058
$ws.Cells.Item(5,9).Value = 0.4068730304461968
$ws.Cells.Item(5,10).Value = 0.4068730304461968
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 14.85604233333333
$ws.Cells.Item(5,14).Value = 44.568127
$ws.Cells.Item(5,15).Value = 0.09286934904108346
$ws.Cells.Item(5,16).Value = 0.09286934904108346
$ws.Cells.Item(5,17).Value = 330.9384382482629
$ws.Cells.Item(5,18).Value = 2978.445944234366
$ws.Cells.Item(5,19).Value = 0.03778603347991123
$ws.Cells.Item(5,20).Value = 0.03778603347991123

# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Tgfb3"
$ws.Cells.Item(6,3).Value = "Tgfbr3"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 22.27635266666667
$ws.Cells.Item(6,8).Value = 66.829058
$ws.Cells.Item(6,9).Value = 0.4068730304461968
$ws.Cells.Item(6,10).Value = 0.4068730304461968
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 114.2734143333333
$ws.Cells.Item(6,14).Value = 342.820243
$ws.Cells.Item(6,15).Value = 0.7143556381787382
$ws.Cells.Item(6,16).Value = 0.7143556381787382
$ws.Cells.Item(6,17).Value = 2545.594878113455
$ws.Cells.Item(6,18).Value = 22910.3539030211
$ws.Cells.Item(6,19).Value = 0.29065204332211
$ws.Cells.Item(6,20).Value = 0.29065204332211

# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Tgfb3"
$ws.Cells.Item(7,3).Value = "Tgfbr3"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 22.27635266666667
$ws.Cells.Item(7,8).Value = 66.829058
$ws.Cells.Item(7,9).Value = 0.4068730304461968
$ws.Cells.Item(7,10).Value = 0.4068730304461968
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 30.83766366666667
$ws.Cells.Item(7,14).Value = 92.512991
$ws.Cells.Item(7,15).Value = 0.1927750127801784
$ws.Cells.Item(7,16).Value = 0.1927750127801784
$ws.Cells.Item(7,17).Value = 686.9506712547199
$ws.Cells.Item(7,18).Value = 6182.556041292478
$ws.Cells.Item(7,19).Value = 0.07843495364417549
$ws.Cells.Item(7,20).Value = 0.07843495364417549

# Row 8
$ws.Cells.Item(8,1).Value = "ECs"
$ws.Cells.Item(8,2).Value = "Tgfb3"
$ws.Cells.Item(8,3).Value = "Tgfbr3"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 31.831433
$ws.Cells.Item(8,8).Value = 95.494299
$ws.Cells.Item(8,9).Value = 0.5813946206523698
$ws.Cells.Item(8,10).Value = 0.5813946206523697
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 14.85604233333333
$ws.Cells.Item(8,14).Value = 44.568127
$ws.Cells.Item(8,15).Value = 0.09286934904108346
$ws.Cells.Item(8,16).Value = 0.09286934904108346
$ws.Cells.Item(8,17).Value = 472.8891161786636
$ws.Cells.Item(8,18).Value = 4256.002045607973
$ws.Cells.Item(8,19).Value = 0.05399373995597324
$ws.Cells.Item(8,20).Value = 0.05399373995597323

# Row 9
$ws.Cells.Item(9,1).Value = "ECs"
$ws.Cells.Item(9,2).Value = "Tgfb3"
$ws.Cells.Item(9,3).Value = "Tgfbr3"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 31.831433
$ws.Cells.Item(9,8).Value = 95.494299
$ws.Cells.Item(9,9).Value = 0.5813946206523698
$ws.Cells.Item(9,10).Value = 0.5813946206523697
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 114.2734143333333
$ws.Cells.Item(9,14).Value = 342.820243
$ws.Cells.Item(9,15).Value = 0.7143556381787382
$ws.Cells.Item(9,16).Value = 0.7143556381787382
$ws.Cells.Item(9,17).Value = 3637.48653203274
$ws.Cells.Item(9,18).Value = 32737.37878829466
$ws.Cells.Item(9,19).Value = 0.415322525269809
$ws.Cells.Item(9,20).Value = 0.4153225252698089

# Row 10
$ws.Cells.Item(10,1).Value = "ECs"
$ws.Cells.Item(10,2).Value = "Tgfb3"
$ws.Cells.Item(10,3).Value = "Tgfbr3"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 31.831433
$ws.Cells.Item(10,8).Value = 95.494299
$ws.Cells.Item(10,9).Value = 0.5813946206523698
$ws.Cells.Item(10,10).Value = 0.5813946206523697
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 30.83766366666667
$ws.Cells.Item(10,14).Value = 92.512991
$ws.Cells.Item(10,15).Value = 0.1927750127801784
$ws.Cells.Item(10,16).Value = 0.1927750127801784
$ws.Cells.Item(10,17).Value = 981.6070248820344
$ws.Cells.Item(10,18).Value = 8834.463223938308
$ws.Cells.Item(10,19).Value = 0.1120783554265876
$ws.Cells.Item(10,20).Value = 0.1120783554265875

